$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5 (shifts everything below down by one)
$ws.Rows("5:5").Insert()

# Populate the newly inserted row with the new verification entry
$ws.Range("A5").Value = "CompFullAndLight.m"
$ws.Range("B5").Value = "Only a small script taken from the GECKO 3 tutorial, no tests deemed necessary."

# Update the selected/active cell to match the saved view state
$ws.Range("B6").Select()
